$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "El sistema desplegará en móvil la geolocalización y en web el catálogo de productos.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El sistema desplegará el menú de inicio", 2)

$d.Content.Find.Execute(
    "El Proveedor selecciona el icono de perfil.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El Proveedor selecciona perfil.", 2)

$d.Content.Find.Execute(
    "El sistema despliega la interfaz de ¿quieres vender tus productos?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El sistema despliega la interfaz de registro", 2)

$d.Content.Find.Execute(
    "El usuario tiene una cuenta activa en la plataforma y puede acceder a ella usando su correo electrónico y contraseña.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El usuario tiene una cuenta activa en la plataforma y puede acceder a ella usando su rut y contraseña.", 2)

$d.Content.Find.Execute(
    "El correo electrónico ingresado ya está registrado: el sistema muestra un mensaje de error indicando que el correo ya está en uso.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El correo electrónico o rut ingresado ya está registrado: el sistema muestra un mensaje de error indicando que el correo ya está en uso.", 2)

Write-Output "done"
